$d = $word.ActiveDocument

# 1. Add "Node.js, " prefix to the technologies list under "Return back to the scene of the crime."
$d.Content.Find.Execute("React, Firebase, AWS", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Node.js, React, Firebase, AWS", 2)

# 2. Shorten the intro sentence about years spent at the company
$d.Content.Find.Execute("I spent last three years in the company called", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "I spent three years in", 2)

# 3. Add "public" before "investment rounds" in the OM Platform description
$d.Content.Find.Execute("Web Application for companies to raise money via investment rounds.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Web Application for companies to raise money via public investment rounds.", 2)

# 4. Change "React.js" to "React" in the EVE project summary line
$d.Content.Find.Execute("TOPMONKS, 4 months, React.js web application", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "TOPMONKS, 4 months, React web application", 2)
